$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.473.19"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.299.62"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.01%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "537.75"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -1.72%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "132.03"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").Value = "2.298.11"
$ws.Range("E9").Value = "  +0.19%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.100"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -1.44%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "5.49"
$r.Style = "Normal"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -0.46%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "23.81"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "2.709.88"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").Value = "58.421.15"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "2.303.60"
$ws.Range("E18").Value = "  -0.90%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "10.57"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  -2.35%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "315.73"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("E23").Value = "  +0.19%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "63.06"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("E26").Value = "  +0.38%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "7.96"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  -1.17%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "170.93"
$r.Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("E31").Value = "  -0.49%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "1.09"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +1.85%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "5.83"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "0.384"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  -0.01%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "17.88"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "
$ws.Range("E37").Value = "  -0.06%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "1.24"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  -0.71%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "291.56"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -2.46%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "140.74"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("E45").Value = "  -0.55%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.556"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -0.11%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "18.33"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -1.52%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.0211"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -2.09%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "10.95"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  +0.69%  "
